# Add two new year columns (2021, 2022) -> O & P to the poverty-level
# table, mirroring the formatting already used for the neighbouring
# existing "year" columns (format-only copy + value write, so we reuse
# existing cell styles instead of minting ad-hoc ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122   # xlPasteFormats

function Copy-FormatOnly {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Row 3 (thick-bottom separator row just above the header row) -----
Copy-FormatOnly "N3" "O3"
Copy-FormatOnly "N3" "P3"

# --- Row 4 (year header row) -------------------------------------------
Copy-FormatOnly "N4" "O4"
$ws.Range("O4").Value = 2021
Copy-FormatOnly "N4" "P4"
$ws.Range("P4").Value = 2022

# --- Row 5 ---------------------------------------------------------------
Copy-FormatOnly "N5" "O5"
$ws.Range("O5").Value = 6.0337796775071091
Copy-FormatOnly "N5" "P5"
$ws.Range("P5").Value = 5.9676405075953687

# --- Row 6 -----------------------------------------------------------
Copy-FormatOnly "N8" "O6"
$ws.Range("O6").Value = 7.3075058743442511
Copy-FormatOnly "N8" "P6"
$ws.Range("P6").Value = 5.6044335798150424

# --- Row 7 -----------------------------------------------------------
Copy-FormatOnly "N8" "O7"
$ws.Range("O7").Value = 5.2767607763499562
Copy-FormatOnly "N8" "P7"
$ws.Range("P7").Value = 6.1789553077823856

# --- Row 8 -----------------------------------------------------------
Copy-FormatOnly "N8" "O8"
$ws.Range("O8").Value = 10.064200140319592
Copy-FormatOnly "N8" "P8"
$ws.Range("P8").Value = 16.5

# --- Row 9 -----------------------------------------------------------
Copy-FormatOnly "N8" "O9"
$ws.Range("O9").Value = 7.5445007460298559
Copy-FormatOnly "N8" "P9"
$ws.Range("P9").Value = 9.1

# --- Row 10 ----------------------------------------------------------
Copy-FormatOnly "N8" "O10"
$ws.Range("O10").Value = 7.9562092224762884
Copy-FormatOnly "N8" "P10"
$ws.Range("P10").Value = 8.8000000000000007

# --- Row 11 ----------------------------------------------------------
Copy-FormatOnly "N8" "O11"
$ws.Range("O11").Value = 8.1696953402867685
Copy-FormatOnly "N8" "P11"
$ws.Range("P11").Value = 6.7

# --- Row 12 ----------------------------------------------------------
Copy-FormatOnly "N8" "O12"
$ws.Range("O12").Value = 2.0701729813092102
Copy-FormatOnly "N8" "P12"
$ws.Range("P12").Value = 0.5

# --- Row 13 ----------------------------------------------------------
Copy-FormatOnly "N8" "O13"
$ws.Range("O13").Value = 2.6482523478927704
Copy-FormatOnly "N8" "P13"
$ws.Range("P13").Value = 2.2000000000000002

# --- Row 14 ----------------------------------------------------------
Copy-FormatOnly "N8" "O14"
$ws.Range("O14").Value = 3.9561647100749857
Copy-FormatOnly "N8" "P14"
$ws.Range("P14").Value = 5.0999999999999996

# --- Row 15 ----------------------------------------------------------
Copy-FormatOnly "N8" "O15"
$ws.Range("O15").Value = 9.4645167179465837
Copy-FormatOnly "N8" "P15"
$ws.Range("P15").Value = 3.9

# --- Row 16 (thick-bottom last data row) ------------------------------
Copy-FormatOnly "N16" "O16"
$ws.Range("O16").Value = 3.1019579996103404
Copy-FormatOnly "N16" "P16"
$ws.Range("P16").Value = 7

# Match the saved selection / active cell from the source workbook.
$ws.Range("P4").Select() | Out-Null

Write-Output "done"
